$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52 (pushes existing rows 52-77 down to 53-78)
$ws.Rows.Item(52).Insert()

# Populate the newly-inserted row 52 with the new weekly price record
$ws.Cells.Item(52, 1).Value  = 10
$ws.Cells.Item(52, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(52, 3).Value  = "La Araucanía"
$ws.Cells.Item(52, 4).Value  = 45027
$ws.Cells.Item(52, 5).Value  = 9
$ws.Cells.Item(52, 6).Value  = 100112010
$ws.Cells.Item(52, 7).Value  = "Achicoria"
$ws.Cells.Item(52, 8).Value  = "Sin especificar"
$ws.Cells.Item(52, 9).Value  = "Primera"
$ws.Cells.Item(52, 10).Value = 100
$ws.Cells.Item(52, 11).Value = 10000
$ws.Cells.Item(52, 12).Value = 10000
$ws.Cells.Item(52, 13).Value = 10000
$ws.Cells.Item(52, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(52, 15).Value = "Región Metropolitana"
$ws.Cells.Item(52, 16).Value = 556
$ws.Cells.Item(52, 17).Value = 18
$ws.Cells.Item(52, 18).Value = "Hortaliza"
